$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7: set the Actual Start-date (E7) to 2019-09-23 (serial 43731),
# matching the date formatting already used by D7.
$ws.Range("E7").Value = 43731
$ws.Range("E7").NumberFormat = $ws.Range("D7").NumberFormat

# Row 7: Status (G7) moves from "Not started" to "In-Progress", shown in
# green text (RGB 0,176,80 = 5287936).
$ws.Range("G7").Value = "In-Progress"
$ws.Range("G7").Font.Color = 5287936

# Move the active selection to F20.
$ws.Range("F20").Select() | Out-Null
